$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C35").Value = "High"
$ws.Range("D35").Value = "Pass"
$ws.Range("C37").Value = "High"
$ws.Range("D37").Value = "Pass"
$ws.Range("C68").Value = "High"
$ws.Range("D68").Value = "Pass"
$ws.Range("D72").Value = "Pass"
